$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Student record update (row 2): Meena Devi / 4-A -> Dhurvan Baskar / 1-B ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Dhurvan"
$ws.Range("C2").Value = "Baskar"
$ws.Range("D2").Value = "1-B"

# Center-align the new data row (creates the new cellXfs entry used by A2:D2)
$ws.Range("A2:D2").HorizontalAlignment = -4108

# --- Column widths (best-fit sizing approximating the authored layout) ---
$ws.Range("A1").ColumnWidth = 2
$ws.Range("B1").ColumnWidth = 9
$ws.Range("C1").ColumnWidth = 9
$ws.Range("D1").ColumnWidth = 4.3333333333333333

# --- Selection moved to C5 ---
$null = $ws.Range("C5").Select()
